$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timing log rows appended below the existing data (rows 2-3 -> now 2-6).
$rows = @(
    @("2026-01-09T13:27:54.321641+00:00", "GET", "/api/inspiration-boards/boards/6960eb294ad6a4df36746c2b", "Inspiration Boards", 0.1207, 200, "6925b1e37b5978266363464e"),
    @("2026-01-09T13:27:58.901902+00:00", "GET", "/api/inspiration-boards/boards", "Inspiration Boards", 0.4167, 200, "6925b1e37b5978266363464e"),
    @("2026-01-09T13:27:59.273055+00:00", "GET", "/api/inspiration-boards/boards", "Inspiration Boards", 0.3671, 200, "6925b1e37b5978266363464e")
)

$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $ws.Cells.Item($r, 7).Value = $rowData[6]

    # "error" column (H) is an empty-text cell (not a blank cell) in the
    # source data, matching rows 2 and 3 above it. A bare "" assignment
    # collapses to a truly blank cell, so force text via a quote-prefixed
    # empty value, then reset the style so no stray quotePrefix formatting
    # is left behind on the cell.
    $ws.Cells.Item($r, 8).Value = "'"
    $ws.Cells.Item($r, 8).Style = "Normal"
}
